$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Q0)
$ws.Range("B2").Value = 0.5440079918808879
$ws.Range("C2").Value = 0.5440079918808879
$ws.Range("D2").Value = 0.3686258812966464
$ws.Range("E2").Value = 0.6071456837503223
$ws.Range("F2").Value = 0.2797714183633027
$ws.Range("G2").Value = 14

# Row 3 (Q1)
$ws.Range("B3").Value = 0.371620442828777
$ws.Range("C3").Value = 0.3874993246909071
$ws.Range("D3").Value = 0.1981274257140365
$ws.Range("E3").Value = 0.4451150701942549
$ws.Range("F3").Value = 0.2550055127533425
$ws.Range("G3").Value = 13

# Row 4 (Q2)
$ws.Range("B4").Value = 0.3184154172361096
$ws.Range("C4").Value = 0.3634107558758037
$ws.Range("D4").Value = 0.1745848588293675
$ws.Range("E4").Value = 0.4178335300444035
$ws.Range("F4").Value = 0.2825786729951437
$ws.Range("G4").Value = 12

# Row 5 (Q3)
$ws.Range("B5").Value = 0.3975283727593981
$ws.Range("C5").Value = 0.4145492339343189
$ws.Range("D5").Value = 0.2091250642104284
$ws.Range("E5").Value = 0.4573019398717093
$ws.Range("F5").Value = 0.2370777989771771
$ws.Range("G5").Value = 11

# Row 6 (Q4)
$ws.Range("B6").Value = 0.3560805627268342
$ws.Range("C6").Value = 0.3805953001122874
$ws.Range("D6").Value = 0.1781696942431756
$ws.Range("E6").Value = 0.4221015212519088
$ws.Range("F6").Value = 0.2389242722689363
$ws.Range("G6").Value = 10

# Row 7 (Q5)
$ws.Range("B7").Value = 0.3392710727664942
$ws.Range("C7").Value = 0.3699241520168967
$ws.Range("D7").Value = 0.1693631219851699
$ws.Range("E7").Value = 0.4115375098155329
$ws.Range("F7").Value = 0.2470638456253208
$ws.Range("G7").Value = 9

# Row 8 (Q6)
$ws.Range("B8").Value = 0.352534700000344
$ws.Range("C8").Value = 0.3776340216059933
$ws.Range("D8").Value = 0.1795390275838225
$ws.Range("E8").Value = 0.4237204592462139
$ws.Range("F8").Value = 0.251301328242729
$ws.Range("G8").Value = 8

# Row 9 (Q7)
$ws.Range("B9").Value = 0.323074747095869
$ws.Range("C9").Value = 0.3439312876808158
$ws.Range("D9").Value = 0.1492574881197667
$ws.Range("E9").Value = 0.3863385666999435
$ws.Range("F9").Value = 0.2288235751843726
$ws.Range("G9").Value = 7

# Row 10 (Q8)
$ws.Range("B10").Value = 0.3581453753460085
$ws.Range("C10").Value = 0.3775233390900978
$ws.Range("D10").Value = 0.1783793918318441
$ws.Range("E10").Value = 0.4223498453081807
$ws.Range("F10").Value = 0.2452214067738235
$ws.Range("G10").Value = 6

# Row 11 (Q9)
$ws.Range("B11").Value = 0.3599616089605324
$ws.Range("C11").Value = 0.3766256120414225
$ws.Range("D11").Value = 0.1882494404226522
$ws.Range("E11").Value = 0.4338772181420132
$ws.Range("F11").Value = 0.2708253138491602
$ws.Range("G11").Value = 5
